$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("138").Insert()

$ws.Cells.Item(138, 1).Value = 10
$ws.Cells.Item(138, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(138, 3).Value = "La Araucanía"
$ws.Cells.Item(138, 4).Value = 44825
$ws.Cells.Item(138, 5).Value = 9
$ws.Cells.Item(138, 6).Value = 100112012
$ws.Cells.Item(138, 7).Value = "Espinaca"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 30
$ws.Cells.Item(138, 11).Value = 10000
$ws.Cells.Item(138, 12).Value = 10000
$ws.Cells.Item(138, 13).Value = 10000
$ws.Cells.Item(138, 14).Value = "$/docena de atados"
$ws.Cells.Item(138, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(138, 16).Value = 3333
$ws.Cells.Item(138, 17).Value = 3
$ws.Cells.Item(138, 18).Value = "Hortaliza"
